# Applies the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions".
#
# For each changed cell we assign the literal text that should appear in the
# sheet. Several Price values (column D) are plain decimal-looking strings
# (e.g. "1.00", "8.30", "0.0000112") that the Excel object model would
# otherwise silently coerce into real numbers (dropping trailing zeros /
# switching to scientific notation). To keep them as text - matching the
# original workbook, where every data cell is stored as a string - we
# temporarily force those specific cells to a text number format before
# assigning the value, then restore the cell to the unstyled "Normal" style
# so the saved file does not pick up a stray formatting difference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

function Set-NumericLookingText($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-PlainText "D2" "50.087.30"
Set-PlainText "E2" "  -17.61%  "
Set-PlainText "D3" "2.232.14"
Set-PlainText "E3" "  -23.13%  "
Set-PlainText "E4" "  +0.17%  "
Set-NumericLookingText "D5" "415.08"
Set-PlainText "E5" "  -21.19%  "
Set-NumericLookingText "D6" "114.56"
Set-PlainText "E6" "  -20.35%  "
Set-PlainText "E7" "  +0.11%  "
Set-NumericLookingText "D8" "0.443"
Set-PlainText "E8" "  -18.97%  "
Set-PlainText "D9" "2.234.09"
Set-PlainText "E9" "  -23.28%  "
Set-NumericLookingText "D10" "4.93"
Set-PlainText "E10" "  -19.96%  "
Set-NumericLookingText "D11" "0.0837"
Set-PlainText "E11" "  -21.95%  "
Set-NumericLookingText "D12" "0.289"
Set-PlainText "E12" "  -19.21%  "
Set-NumericLookingText "D13" "0.119"
Set-PlainText "E13" "  -7.55%  "
Set-PlainText "D14" "2.616.24"
Set-PlainText "E14" "  -23.33%  "
Set-PlainText "D15" "50.404.07"
Set-PlainText "E15" "  -17.15%  "
Set-NumericLookingText "D16" "17.82"
Set-PlainText "E16" "  -20.89%  "
Set-PlainText "B17" "ShibaInu"
Set-PlainText "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-NumericLookingText "D17" "0.0000112"
Set-PlainText "E17" "  -20.87%  "
Set-PlainText "B18" "WrappedEther"
Set-PlainText "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-PlainText "D18" "2.225.75"
Set-PlainText "E18" "  -23.46%  "
Set-NumericLookingText "D19" "3.86"
Set-PlainText "E19" "  -21.21%  "
Set-NumericLookingText "D20" "283.73"
Set-PlainText "E20" "  -19.64%  "
Set-NumericLookingText "D21" "1.00"
Set-PlainText "E21" "  +0.26%  "
Set-NumericLookingText "D22" "5.63"
Set-PlainText "E22" "  -1.65%  "
Set-NumericLookingText "D23" "8.30"
Set-PlainText "E23" "  -28.07%  "
Set-NumericLookingText "D24" "4.90"
Set-PlainText "E24" "  -24.80%  "
Set-NumericLookingText "D25" "0.998"
Set-PlainText "E25" "  -0.09%  "
Set-NumericLookingText "D26" "51.35"
Set-PlainText "E26" "  -20.74%  "
Set-NumericLookingText "D27" "0.355"
Set-PlainText "E27" "  -21.11%  "
Set-PlainText "D28" "2.321.37"
Set-PlainText "E28" "  -23.53%  "
Set-NumericLookingText "D29" "0.997"
Set-PlainText "E29" "  -0.21%  "
Set-NumericLookingText "D30" "0.128"
Set-PlainText "E30" "  -28.10%  "
Set-NumericLookingText "D31" "6.50"
Set-PlainText "E31" "  -17.12%  "
Set-PlainText "D32" "0.0₃0619"
Set-PlainText "E32" "  -28.38%  "
Set-NumericLookingText "D33" "140.04"
Set-PlainText "E33" "  -8.35%  "
Set-NumericLookingText "D34" "16.11"
Set-PlainText "E34" "  -17.94%  "
Set-NumericLookingText "D35" "1.27"
Set-PlainText "E35" "  -24.67%  "
Set-PlainText "E36" "  -18.93%  "
Set-NumericLookingText "D37" "0.999"
Set-PlainText "E37" "  +0.16%  "
Set-PlainText "E38" "  -27.70%  "
Set-NumericLookingText "D39" "31.58"
Set-PlainText "E39" "  -15.84%  "
Set-PlainText "B40" "ImmutableX"
Set-PlainText "C40" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-NumericLookingText "D40" "0.934"
Set-PlainText "E40" "  -21.75%  "
Set-PlainText "B41" "Fetch.AI"
Set-PlainText "C41" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-NumericLookingText "D41" "0.723"
Set-PlainText "E41" "  -27.16%  "
Set-NumericLookingText "D42" "10.10"
Set-PlainText "E42" "  -2.30%  "
Set-PlainText "E43" "  -17.95%  "
Set-NumericLookingText "D44" "3.00"
Set-PlainText "E44" "  -18.74%  "
Set-NumericLookingText "D45" "0.0477"
Set-PlainText "E45" "  -17.92%  "
Set-PlainText "D46" "1.811.37"
Set-PlainText "E46" "  -20.47%  "
Set-NumericLookingText "D47" "1.08"
Set-PlainText "E47" "  -25.98%  "
Set-PlainText "E48" "  -18.31%  "
Set-NumericLookingText "D49" "0.0773"
Set-PlainText "E49" "  -15.69%  "
Set-NumericLookingText "D50" "4.60"
Set-PlainText "E50" "  -6.37%  "
Set-NumericLookingText "D51" "15.03"
Set-PlainText "E51" "  -25.72%  "
